$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1681.8108
$ws.Range("I15").Value = 1681.8108
$ws.Range("K15").Value = 5045.4324
$ws.Range("M15").Value = -4876.4324
$ws.Range("H17").Value = 65837.42
$ws.Range("J17").Value = 65837.42
$ws.Range("L17").Value = 197512.26
$ws.Range("N17").Value = -197848.26
$ws.Range("H34").Value = 15461.5
$ws.Range("I34").Value = 15461.5
$ws.Range("K34").Value = 15461.5
$ws.Range("M34").Value = -15258.5
$ws.Range("H36").Value = 15461.5
$ws.Range("I36").Value = 15461.5
$ws.Range("K36").Value = 15461.5
$ws.Range("M36").Value = -14746.5
$ws.Range("H38").Value = 1019.3333
$ws.Range("I38").Value = 1019.3333
$ws.Range("K38").Value = 3057.9999
$ws.Range("M38").Value = -2685.9999
$ws.Range("H69").Value = 23925
$ws.Range("J69").Value = 25900
$ws.Range("L69").Value = 77700
$ws.Range("N69").Value = -79448
$ws.Range("H72").Value = 23925
$ws.Range("J72").Value = 25900
$ws.Range("L72").Value = 233100
$ws.Range("N72").Value = -241836
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 77773.5
$ws.Range("J126").Value = 77773.5
$ws.Range("L126").Value = 77773.5
$ws.Range("N126").Value = -87653.5
$ws.Range("H130").Value = 83076.92
$ws.Range("J130").Value = 83076.92
$ws.Range("L130").Value = 83076.92
$ws.Range("N130").Value = -93116.92
$ws.Range("H133").Value = 49266.25
$ws.Range("J133").Value = 49266.25
$ws.Range("L133").Value = 49266.25
$ws.Range("N133").Value = -59386.25
$ws.Range("H134").Value = 72166.61
$ws.Range("J134").Value = 72166.61
$ws.Range("L134").Value = 72166.61
$ws.Range("N134").Value = -82306.61
$ws.Range("H136").Value = 29999
$ws.Range("I136").Value = 29999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 29999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -24899
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 7166.808
$ws.Range("I137").Value = 1207.2646
$ws.Range("J137").Value = 12362.308
$ws.Range("K137").Value = 3621.7938
$ws.Range("L137").Value = 37086.924
$ws.Range("M137").Value = -1071.7938
$ws.Range("N137").Value = -42186.924
$ws.Range("H138").Value = 4719.7163
$ws.Range("I138").Value = 1484.1818
$ws.Range("J138").Value = 6301.533
$ws.Range("K138").Value = 4452.5454
$ws.Range("L138").Value = 18904.599
$ws.Range("M138").Value = 687.4546
$ws.Range("N138").Value = -29184.599
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 2747.6924
$ws.Range("J141").Value = 2046
$ws.Range("L141").Value = 6138
$ws.Range("N141").Value = -16498

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3452.853
$ws.Range("I32").Value = 3452.853
$ws.Range("K32").Value = 3452.853
$ws.Range("M32").Value = -3165.853
$ws.Range("H74").Value = 168228.4
$ws.Range("I74").Value = 186266.33
$ws.Range("K74").Value = 186266.33
$ws.Range("M74").Value = -185392.33
$ws.Range("H77").Value = 168228.4
$ws.Range("I77").Value = 186266.33
$ws.Range("K77").Value = 931331.6499999999
$ws.Range("M77").Value = -926963.6499999999
$ws.Range("H80").Value = 79526.42999999999
$ws.Range("J80").Value = 79526.42999999999
$ws.Range("L80").Value = 79526.42999999999
$ws.Range("N80").Value = -81522.42999999999
$ws.Range("H83").Value = 79526.42999999999
$ws.Range("J83").Value = 79526.42999999999
$ws.Range("L83").Value = 238579.29
$ws.Range("N83").Value = -248563.29
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H140").Value = 43867.8
$ws.Range("J140").Value = 43867.8
$ws.Range("L140").Value = 43867.8
$ws.Range("N140").Value = -54227.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 27123.75
$ws.Range("I134").Value = 1279.6765
$ws.Range("J134").Value = 114993.6
$ws.Range("K134").Value = 3839.0295
$ws.Range("L134").Value = 344980.8
$ws.Range("M134").Value = -1304.0295
$ws.Range("N134").Value = -350050.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 17885
$ws.Range("J28").Value = 17885
$ws.Range("L28").Value = 17885
$ws.Range("N28").Value = -18375
$ws.Range("H31").Value = 458791
$ws.Range("I31").Value = 5650680.5
$ws.Range("K31").Value = 5650680.5
$ws.Range("M31").Value = -5650385.5
$ws.Range("H34").Value = 458791
$ws.Range("I34").Value = 5650680.5
$ws.Range("K34").Value = 5650680.5
$ws.Range("M34").Value = -5650478.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 650.4167
$ws.Range("I92").Value = 457.75
$ws.Range("J92").Value = 746.75
$ws.Range("K92").Value = 1373.25
$ws.Range("L92").Value = 2240.25
$ws.Range("M92").Value = -125.25
$ws.Range("N92").Value = -4736.25
$ws.Range("H131").Value = 41039.17
$ws.Range("I131").Value = 78165.16
$ws.Range("J131").Value = 28663.846
$ws.Range("K131").Value = 234495.48
$ws.Range("L131").Value = 85991.538
$ws.Range("M131").Value = -229455.48
$ws.Range("N131").Value = -96071.538

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 700485.9399999999
$ws.Range("I113").Value = 1309038.2
$ws.Range("K113").Value = 1309038.2
$ws.Range("M113").Value = -1306868.2
$ws.Range("H122").Value = 413041.25
$ws.Range("I122").Value = 505050.62
$ws.Range("K122").Value = 1515151.86
$ws.Range("M122").Value = -1512701.86

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3311.8823
$ws.Range("I46").Value = 3144.611
$ws.Range("K46").Value = 3144.611
$ws.Range("M46").Value = -2956.611
$ws.Range("H61").Value = 6135.7188
$ws.Range("I61").Value = 5386.231
$ws.Range("K61").Value = 5386.231
$ws.Range("M61").Value = -5184.231
$ws.Range("H113").Value = 6135.7188
$ws.Range("I113").Value = 5386.231
$ws.Range("K113").Value = 5386.231
$ws.Range("M113").Value = -3216.231

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 847.5
$ws.Range("J8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("N8").Value = -1080
$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 20000
$ws.Range("K39").Value = 20000
$ws.Range("M39").Value = -19587
$ws.Range("H74").Value = 10990.625
$ws.Range("I74").Value = 9990
$ws.Range("J74").Value = 11324.167
$ws.Range("K74").Value = 9990
$ws.Range("L74").Value = 11324.167
$ws.Range("M74").Value = -9054
$ws.Range("N74").Value = -13196.167
$ws.Range("H77").Value = 10990.625
$ws.Range("I77").Value = 9990
$ws.Range("J77").Value = 11324.167
$ws.Range("K77").Value = 29970
$ws.Range("L77").Value = 33972.501
$ws.Range("M77").Value = -25290
$ws.Range("N77").Value = -43332.501
$ws.Range("H122").Value = 4137.9487
$ws.Range("I122").Value = 3530.6562
$ws.Range("K122").Value = 10591.9686
$ws.Range("M122").Value = -8141.9686
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value = 58000
$ws.Range("J141").Value = 58000
$ws.Range("L141").Value = 58000
$ws.Range("N141").Value = -68360
